$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting the existing table (old A:AC) to (B:AD)
$ws.Columns("A:A").Insert()

# New column header (row 3 holds the "pandas column name" header row) -- bold, no border
$ws.Range("A3").Value = "Match ID"
$ws.Range("A3").Font.Bold = $true

# Data rows 4-20: visible player rows, all belong to Match ID 24
$ws.Range("A4:A20").Value = 24
$ws.Range("A4:A20").Font.Bold = $true

# Row 21 is the hidden "totals" row; its new Match ID cell keeps the plain/default style
$ws.Range("A21").Value = 24

# Restore natural row heights (insert + value writes can otherwise pin an explicit
# height on the touched hidden rows); this brings rows back to sheet default sizing.
$ws.Rows("1:21").AutoFit()

# Update the selection to match the post-edit state
$ws.Range("A3:A20").Select() | Out-Null
